$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Update the "Estado de Cuenta" data rows (new part of the account-statement
# database): rows 17 and 19 swap their "Periodo Mora" / "Valor Mora" values,
# row 18 stays as-is.
$ws.Range("E17").Value = "1809"
$ws.Range("F17").Value = 31249

$ws.Range("E19").Value = "1902"
$ws.Range("F19").Value = 28124

# --- Nudge the logo image 19pt to the left (its size and vertical position
# are unchanged). Re-assert the width explicitly afterwards using its exact
# point value (975600 EMU / 612000 EMU at 12700 EMU-per-point): this COM
# layer recomputes the shape's EMU extent whenever a geometry setter runs,
# and the recompute drifts unless Width is pinned back to its true value
# (the rounded read-back from $shp.Width is not precise enough to avoid it).
$shp = $ws.Shapes.Item(1)
$shp.Left = $shp.Left - 19
$shp.Width = 76.81889763779527

